# contratos-7-2014.xlsx - "fix: fixed formatting when scrapping floating point numbers"
#
# Two independent clean-ups of text that was scraped from the source site:
#
#   1) A handful of "Razon social" / "Nombre Fantasia" cells (column E / F)
#      used a comma to separate co-contracted people/companies, which reads
#      ambiguously next to "Apellido, Nombre" formatting. Those commas are
#      normalised to periods.
#
#   2) The "Importe" column (H) was scraped as text using Spanish/Argentine
#      number formatting (thousands separator "." , decimal separator ",").
#      That is replaced with a plain decimal representation: no thousands
#      separator, "." as the decimal separator (e.g. "2.400,00" -> "2400.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Proveedor name fields: comma used as a separator -> period ---------
$nameFixes = @(
    @{ Row = 108; Col = 5; Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO" },
    @{ Row = 143; Col = 5; Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA" },
    @{ Row = 143; Col = 6; Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA" },
    @{ Row = 153; Col = 5; Value = "RICCOTTI. MARIANA EDITH" },
    @{ Row = 199; Col = 5; Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN" },
    @{ Row = 212; Col = 5; Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH" },
    @{ Row = 164; Col = 6; Value = "MERCANZINI. GASTON ARIEL" }
)

foreach ($fix in $nameFixes) {
    $ws.Cells.Item($fix.Row, $fix.Col).Value = $fix.Value
}

# --- 2) "Importe" column (H, rows 2-252): re-format the numbers ------------
# The values are stored as text (not numbers) in the sheet. A plain
# `.Value = "2400.00"` assignment would let Excel interpret the numeric-
# looking text and coerce the cell into a real number (dropping trailing
# zeros / the intended text formatting). To keep these as text, the range is
# temporarily marked as Text ("@") while the values are rewritten, then the
# style is restored to Normal so cell formatting matches the original file.
$importeRange = $ws.Range("H2:H252")
$importeRange.NumberFormat = "@"

for ($row = 2; $row -le 252; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $oldText = $cell.Text
    # "2.680.000,00" -> remove thousands "." -> "2680000,00" -> "," becomes
    # the decimal point -> "2680000.00"
    $newText = $oldText.Replace(".", "").Replace(",", ".")
    $cell.Value = $newText
}

$importeRange.Style = "Normal"
